$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2:E24").Value = "fullRNASeq"
